$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update weekly price data for rows 63-80 (Espárragos, Femacal de La Calera).
# Existing rows 63-79 shift to new values; row 80 is a newly appended record.

# Row 63
$ws.Range("A63").Value = 3
$ws.Range("B63").Value = "Femacal de La Calera"
$ws.Range("C63").Value = "Coquimbo"
$ws.Range("D63").Value = 45215
$ws.Range("E63").Value = 5
$ws.Range("F63").Value = 300000000
$ws.Range("G63").Value = "Espárragos"
$ws.Range("H63").Value = "Verde"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 1200
$ws.Range("K63").Value = 1800
$ws.Range("L63").Value = 1800
$ws.Range("M63").Value = 1800
$ws.Range("N63").Value = "`$/kilo"
$ws.Range("O63").Value = "Provincia de Quillota"
$ws.Range("P63").Value = 1800
$ws.Range("Q63").Value = 1
$ws.Range("R63").Value = "Hortaliza"

# Row 64
$ws.Range("A64").Value = 3
$ws.Range("B64").Value = "Femacal de La Calera"
$ws.Range("C64").Value = "Coquimbo"
$ws.Range("D64").Value = 44923
$ws.Range("E64").Value = 5
$ws.Range("F64").Value = 300000000
$ws.Range("G64").Value = "Espárragos"
$ws.Range("H64").Value = "Verde"
$ws.Range("I64").Value = "Primera"
$ws.Range("J64").Value = 480
$ws.Range("K64").Value = 1500
$ws.Range("L64").Value = 1500
$ws.Range("M64").Value = 1500
$ws.Range("N64").Value = "`$/kilo"
$ws.Range("O64").Value = "Provincia de Quillota"
$ws.Range("P64").Value = 1500
$ws.Range("Q64").Value = 1
$ws.Range("R64").Value = "Hortaliza"

# Row 65
$ws.Range("A65").Value = 3
$ws.Range("B65").Value = "Femacal de La Calera"
$ws.Range("C65").Value = "Coquimbo"
$ws.Range("D65").Value = 44923
$ws.Range("E65").Value = 5
$ws.Range("F65").Value = 300000000
$ws.Range("G65").Value = "Espárragos"
$ws.Range("H65").Value = "Verde"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 480
$ws.Range("K65").Value = 1500
$ws.Range("L65").Value = 1500
$ws.Range("M65").Value = 1500
$ws.Range("N65").Value = "`$/kilo"
$ws.Range("O65").Value = "Perú"
$ws.Range("P65").Value = 1500
$ws.Range("Q65").Value = 1
$ws.Range("R65").Value = "Hortaliza"

# Row 66
$ws.Range("A66").Value = 3
$ws.Range("B66").Value = "Femacal de La Calera"
$ws.Range("C66").Value = "Coquimbo"
$ws.Range("D66").Value = 44914
$ws.Range("E66").Value = 5
$ws.Range("F66").Value = 300000000
$ws.Range("G66").Value = "Espárragos"
$ws.Range("H66").Value = "Verde"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 680
$ws.Range("K66").Value = 1500
$ws.Range("L66").Value = 1500
$ws.Range("M66").Value = 1500
$ws.Range("N66").Value = "`$/kilo"
$ws.Range("O66").Value = "Perú"
$ws.Range("P66").Value = 1500
$ws.Range("Q66").Value = 1
$ws.Range("R66").Value = "Hortaliza"

# Row 67
$ws.Range("A67").Value = 3
$ws.Range("B67").Value = "Femacal de La Calera"
$ws.Range("C67").Value = "Coquimbo"
$ws.Range("D67").Value = 44841
$ws.Range("E67").Value = 5
$ws.Range("F67").Value = 300000000
$ws.Range("G67").Value = "Espárragos"
$ws.Range("H67").Value = "Verde"
$ws.Range("I67").Value = "Primera"
$ws.Range("J67").Value = 1260
$ws.Range("K67").Value = 1400
$ws.Range("L67").Value = 1500
$ws.Range("M67").Value = 1454
$ws.Range("N67").Value = "`$/kilo"
$ws.Range("O67").Value = "Provincia de Quillota"
$ws.Range("P67").Value = 1454
$ws.Range("Q67").Value = 1
$ws.Range("R67").Value = "Hortaliza"

# Row 68
$ws.Range("A68").Value = 3
$ws.Range("B68").Value = "Femacal de La Calera"
$ws.Range("C68").Value = "Coquimbo"
$ws.Range("D68").Value = 45204
$ws.Range("E68").Value = 5
$ws.Range("F68").Value = 300000000
$ws.Range("G68").Value = "Espárragos"
$ws.Range("H68").Value = "Verde"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 1200
$ws.Range("K68").Value = 1800
$ws.Range("L68").Value = 1800
$ws.Range("M68").Value = 1800
$ws.Range("N68").Value = "`$/kilo"
$ws.Range("O68").Value = "Provincia de Quillota"
$ws.Range("P68").Value = 1800
$ws.Range("Q68").Value = 1
$ws.Range("R68").Value = "Hortaliza"

# Row 69
$ws.Range("A69").Value = 3
$ws.Range("B69").Value = "Femacal de La Calera"
$ws.Range("C69").Value = "Coquimbo"
$ws.Range("D69").Value = 45195
$ws.Range("E69").Value = 5
$ws.Range("F69").Value = 300000000
$ws.Range("G69").Value = "Espárragos"
$ws.Range("H69").Value = "Verde"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 1300
$ws.Range("K69").Value = 2000
$ws.Range("L69").Value = 2000
$ws.Range("M69").Value = 2000
$ws.Range("N69").Value = "`$/kilo"
$ws.Range("O69").Value = "Provincia de Linares"
$ws.Range("P69").Value = 2000
$ws.Range("Q69").Value = 1
$ws.Range("R69").Value = "Hortaliza"

# Row 70
$ws.Range("A70").Value = 3
$ws.Range("B70").Value = "Femacal de La Calera"
$ws.Range("C70").Value = "Coquimbo"
$ws.Range("D70").Value = 45211
$ws.Range("E70").Value = 5
$ws.Range("F70").Value = 300000000
$ws.Range("G70").Value = "Espárragos"
$ws.Range("H70").Value = "Verde"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 1200
$ws.Range("K70").Value = 1800
$ws.Range("L70").Value = 1800
$ws.Range("M70").Value = 1800
$ws.Range("N70").Value = "`$/kilo"
$ws.Range("O70").Value = "Provincia de Linares"
$ws.Range("P70").Value = 1800
$ws.Range("Q70").Value = 1
$ws.Range("R70").Value = "Hortaliza"

# Row 71
$ws.Range("A71").Value = 3
$ws.Range("B71").Value = "Femacal de La Calera"
$ws.Range("C71").Value = "Coquimbo"
$ws.Range("D71").Value = 44838
$ws.Range("E71").Value = 5
$ws.Range("F71").Value = 300000000
$ws.Range("G71").Value = "Espárragos"
$ws.Range("H71").Value = "Verde"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 1480
$ws.Range("K71").Value = 1400
$ws.Range("L71").Value = 1500
$ws.Range("M71").Value = 1461
$ws.Range("N71").Value = "`$/kilo"
$ws.Range("O71").Value = "Provincia de Quillota"
$ws.Range("P71").Value = 1461
$ws.Range("Q71").Value = 1
$ws.Range("R71").Value = "Hortaliza"

# Row 72
$ws.Range("A72").Value = 3
$ws.Range("B72").Value = "Femacal de La Calera"
$ws.Range("C72").Value = "Coquimbo"
$ws.Range("D72").Value = 45209
$ws.Range("E72").Value = 5
$ws.Range("F72").Value = 300000000
$ws.Range("G72").Value = "Espárragos"
$ws.Range("H72").Value = "Verde"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 1800
$ws.Range("K72").Value = 1800
$ws.Range("L72").Value = 1800
$ws.Range("M72").Value = 1800
$ws.Range("N72").Value = "`$/kilo"
$ws.Range("O72").Value = "Provincia de Linares"
$ws.Range("P72").Value = 1800
$ws.Range("Q72").Value = 1
$ws.Range("R72").Value = "Hortaliza"

# Row 73
$ws.Range("A73").Value = 3
$ws.Range("B73").Value = "Femacal de La Calera"
$ws.Range("C73").Value = "Coquimbo"
$ws.Range("D73").Value = 44859
$ws.Range("E73").Value = 5
$ws.Range("F73").Value = 300000000
$ws.Range("G73").Value = "Espárragos"
$ws.Range("H73").Value = "Verde"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 1580
$ws.Range("K73").Value = 1400
$ws.Range("L73").Value = 1400
$ws.Range("M73").Value = 1400
$ws.Range("N73").Value = "`$/kilo"
$ws.Range("O73").Value = "Provincia de Linares"
$ws.Range("P73").Value = 1400
$ws.Range("Q73").Value = 1
$ws.Range("R73").Value = "Hortaliza"

# Row 74
$ws.Range("A74").Value = 3
$ws.Range("B74").Value = "Femacal de La Calera"
$ws.Range("C74").Value = "Coquimbo"
$ws.Range("D74").Value = 44845
$ws.Range("E74").Value = 5
$ws.Range("F74").Value = 300000000
$ws.Range("G74").Value = "Espárragos"
$ws.Range("H74").Value = "Verde"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 1000
$ws.Range("K74").Value = 1300
$ws.Range("L74").Value = 1500
$ws.Range("M74").Value = 1396
$ws.Range("N74").Value = "`$/kilo"
$ws.Range("O74").Value = "Provincia de Quillota"
$ws.Range("P74").Value = 1396
$ws.Range("Q74").Value = 1
$ws.Range("R74").Value = "Hortaliza"

# Row 75
$ws.Range("A75").Value = 3
$ws.Range("B75").Value = "Femacal de La Calera"
$ws.Range("C75").Value = "Coquimbo"
$ws.Range("D75").Value = 45212
$ws.Range("E75").Value = 5
$ws.Range("F75").Value = 300000000
$ws.Range("G75").Value = "Espárragos"
$ws.Range("H75").Value = "Verde"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 1100
$ws.Range("K75").Value = 1800
$ws.Range("L75").Value = 1800
$ws.Range("M75").Value = 1800
$ws.Range("N75").Value = "`$/kilo"
$ws.Range("O75").Value = "Provincia de Linares"
$ws.Range("P75").Value = 1800
$ws.Range("Q75").Value = 1
$ws.Range("R75").Value = "Hortaliza"

# Row 76
$ws.Range("A76").Value = 3
$ws.Range("B76").Value = "Femacal de La Calera"
$ws.Range("C76").Value = "Coquimbo"
$ws.Range("D76").Value = 44181
$ws.Range("E76").Value = 5
$ws.Range("F76").Value = 300000000
$ws.Range("G76").Value = "Espárragos"
$ws.Range("H76").Value = "Verde"
$ws.Range("I76").Value = "Primera"
$ws.Range("J76").Value = 1000
$ws.Range("K76").Value = 1300
$ws.Range("L76").Value = 1300
$ws.Range("M76").Value = 1300
$ws.Range("N76").Value = "`$/kilo"
$ws.Range("O76").Value = "Provincia de Quillota"
$ws.Range("P76").Value = 1300
$ws.Range("Q76").Value = 1
$ws.Range("R76").Value = "Hortaliza"

# Row 77
$ws.Range("A77").Value = 3
$ws.Range("B77").Value = "Femacal de La Calera"
$ws.Range("C77").Value = "Coquimbo"
$ws.Range("D77").Value = 44181
$ws.Range("E77").Value = 5
$ws.Range("F77").Value = 300000000
$ws.Range("G77").Value = "Espárragos"
$ws.Range("H77").Value = "Verde"
$ws.Range("I77").Value = "Segunda"
$ws.Range("J77").Value = 900
$ws.Range("K77").Value = 900
$ws.Range("L77").Value = 900
$ws.Range("M77").Value = 900
$ws.Range("N77").Value = "`$/kilo"
$ws.Range("O77").Value = "Provincia de Quillota"
$ws.Range("P77").Value = 900
$ws.Range("Q77").Value = 1
$ws.Range("R77").Value = "Hortaliza"

# Row 78
$ws.Range("A78").Value = 3
$ws.Range("B78").Value = "Femacal de La Calera"
$ws.Range("C78").Value = "Coquimbo"
$ws.Range("D78").Value = 44918
$ws.Range("E78").Value = 5
$ws.Range("F78").Value = 300000000
$ws.Range("G78").Value = "Espárragos"
$ws.Range("H78").Value = "Verde"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 550
$ws.Range("K78").Value = 1500
$ws.Range("L78").Value = 1500
$ws.Range("M78").Value = 1500
$ws.Range("N78").Value = "`$/kilo"
$ws.Range("O78").Value = "Provincia de Quillota"
$ws.Range("P78").Value = 1500
$ws.Range("Q78").Value = 1
$ws.Range("R78").Value = "Hortaliza"

# Row 79
$ws.Range("A79").Value = 3
$ws.Range("B79").Value = "Femacal de La Calera"
$ws.Range("C79").Value = "Coquimbo"
$ws.Range("D79").Value = 45194
$ws.Range("E79").Value = 5
$ws.Range("F79").Value = 300000000
$ws.Range("G79").Value = "Espárragos"
$ws.Range("H79").Value = "Verde"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 1300
$ws.Range("K79").Value = 1900
$ws.Range("L79").Value = 1900
$ws.Range("M79").Value = 1900
$ws.Range("N79").Value = "`$/kilo"
$ws.Range("O79").Value = "Provincia de Linares"
$ws.Range("P79").Value = 1900
$ws.Range("Q79").Value = 1
$ws.Range("R79").Value = "Hortaliza"

# Row 80
$ws.Range("A80").Value = 3
$ws.Range("B80").Value = "Femacal de La Calera"
$ws.Range("C80").Value = "Coquimbo"
$ws.Range("D80").Value = 44900
$ws.Range("E80").Value = 5
$ws.Range("F80").Value = 300000000
$ws.Range("G80").Value = "Espárragos"
$ws.Range("H80").Value = "Verde"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 950
$ws.Range("K80").Value = 1500
$ws.Range("L80").Value = 1500
$ws.Range("M80").Value = 1500
$ws.Range("N80").Value = "`$/kilo"
$ws.Range("O80").Value = "Provincia de Quillota"
$ws.Range("P80").Value = 1500
$ws.Range("Q80").Value = 1
$ws.Range("R80").Value = "Hortaliza"

# Ensure the newly added row 80 date cell uses the same date format as the rest of column D.
$ws.Range("D80").NumberFormat = $ws.Range("D79").NumberFormat

